$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- 1. Stamp the blank template row (currently row 36, all empty) down into
#        rows 46:57 first, while it is still blank, so the new rows inherit
#        the correct number formats / styles (s="6"/"7"/"8") exactly like the
#        pre-existing placeholder rows did. ---
$ws.Range("A36:G36").Copy($ws.Range("A46:G57")) | Out-Null
for ($r = 46; $r -le 57; $r++) {
    $ws.Range("D$r").Formula = "=Tableau2[[#This Row],[Fin]]-Tableau2[[#This Row],[Début]]"
}

# --- 2. Fill in the real journal entries for rows 36-45. Dates/times are
#        written as their underlying serial values so no locale parsing is
#        involved. Column E/F (Titre/Description) text is written in the
#        exact order the shared-string table needs to see each new unique
#        string for the first time, so the resulting sharedStrings.xml
#        ordering lines up. ---

# Row 36
$ws.Range("A36").Value = 44691
$ws.Range("B36").Value = 0.33333333333333331
$ws.Range("C36").Value = 0.39930555555555558

# Row 37
$ws.Range("A37").Value = 44691
$ws.Range("B37").Value = 0.40972222222222227
$ws.Range("C37").Value = 0.44097222222222227

# Row 38
$ws.Range("A38").Value = 44691
$ws.Range("B38").Value = 0.44444444444444442
$ws.Range("C38").Value = 0.47569444444444442

# Row 39
$ws.Range("A39").Value = 44691
$ws.Range("B39").Value = 0.47916666666666669
$ws.Range("C39").Value = 0.51041666666666663

# Row 40
$ws.Range("A40").Value = 44691
$ws.Range("B40").Value = 0.63888888888888895
$ws.Range("C40").Value = 0.67013888888888884

# Row 41
$ws.Range("A41").Value = 44691
$ws.Range("B41").Value = 0.67013888888888884
$ws.Range("C41").Value = 0.70486111111111116

# Row 42
$ws.Range("A42").Value = 44692
$ws.Range("B42").Value = 0.36458333333333331
$ws.Range("C42").Value = 0.39583333333333331

# Row 43
$ws.Range("A43").Value = 44692
$ws.Range("B43").Value = 0.40972222222222227
$ws.Range("C43").Value = 0.45833333333333331

# Row 44
$ws.Range("A44").Value = 44692
$ws.Range("B44").Value = 0.45833333333333331
$ws.Range("C44").Value = 0.47916666666666669

# Row 45
$ws.Range("A45").Value = 44692
$ws.Range("B45").Value = 0.47916666666666669
$ws.Range("C45").Value = 0.51041666666666663

# Titre / Description text, in first-use order.
$ws.Range("E36").Value = "implemenation du template pour gabarit"
$ws.Range("E37").Value = "implemantation de la page d'acceuil"
$ws.Range("E38").Value = "implemantation de la page de login"
$ws.Range("E39").Value = "implementation de la page de signup"
$ws.Range("F36").Value = "avec un peut d'aide de samuel mon collègue de classe"
$ws.Range("E42").Value = "implementation de la fonction login"
$ws.Range("E40").Value = "implementation de la fonction signup"
$ws.Range("E41").Value = "implementation de la fonction signup"
$ws.Range("F40").Value = "j'ai développer une partie de la fonction de signup (récuperer les info, hash password )"
$ws.Range("F41").Value = "j'ai créer les fonction pour ajouter un item dans la base de donnée "
$ws.Range("F42").Value = "création de la fonction tryLogin qui a pour but de se connecter"
$ws.Range("E43").Value = "implementation des helpers"
$ws.Range("E44").Value = "implementation des helpers"
$ws.Range("F43").Value = "dev la fonction helpers qui modifie les boutons de connexion en bouton de deconnexion et de page personnel"
$ws.Range("F44").Value = "dev de la fonction de flashmessage"
$ws.Range("E45").Value = "correction bug affichage "
$ws.Range("F45").Value = "quand je me connecte ca affichais la page d'acceuil et la page de sign up le problème était un break dans l'index"

# --- 3. Rows 40/43/45 wrap onto two lines once the long description text is
#        in place, so they need the taller row height. ---
$ws.Rows.Item(40).RowHeight = 30
$ws.Rows.Item(43).RowHeight = 30
$ws.Rows.Item(45).RowHeight = 30

# --- 4. Grow the table (and its autofilter) to cover the newly added rows. ---
$tbl = $ws.ListObjects.Item("Tableau2")
$tbl.Resize($ws.Range("A1:G57")) | Out-Null

# --- 5. Restore the view state (scroll position + active selection). ---
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 26
$ws.Range("F42").Select() | Out-Null
